# Generate Report for Handoff
# Replace references to the old source file GUID/name with the new one,
# and refresh the associated handoff timestamps, across all three sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "4d15dc09-fd93-4f9e-962a-26db6c0ff672"
$newGuid = "4eb5e790-29b6-4400-9b18-a3684345267b"

$oldZhXlf = "$oldGuid.97d4ba0ec9028312904f710aa1f24d573e9a6246.zh-cn.xlf"
$newZhXlf = "$newGuid.4421b21dab6312a3e2c0564440891d670cbffa4d.zh-cn.xlf"

$oldDeXlf = "$oldGuid.97d4ba0ec9028312904f710aa1f24d573e9a6246.de-de.xlf"
$newDeXlf = "$newGuid.4421b21dab6312a3e2c0564440891d670cbffa4d.de-de.xlf"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-12 07:11:00"

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-08-12 07:10:52"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-08-12 07:11:00"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}
